$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Change the rich-text run "重命名" -> "auto" inside cell B2, while
#    preserving the original run formatting (bold, size, color, font).
# ---------------------------------------------------------------------
$cell = $ws.Range("B2")

$midFont = $cell.Characters(56, 3).Font
$midName = $midFont.Name
$midBold = $midFont.Bold
$midSize = $midFont.Size
$midColor = $midFont.Color

$lastFont = $cell.Characters(59, 2).Font
$lastName = $lastFont.Name
$lastBold = $lastFont.Bold
$lastSize = $lastFont.Size
$lastColor = $lastFont.Color

$cell.Characters(56, 3).Text = "auto"

$newMid = $cell.Characters(56, 4)
$newMid.Font.Name = $midName
$newMid.Font.Bold = $midBold
$newMid.Font.Size = $midSize
$newMid.Font.Color = $midColor

$newLast = $cell.Characters(60, 2)
$newLast.Font.Name = $lastName
$newLast.Font.Bold = $lastBold
$newLast.Font.Size = $lastSize
$newLast.Font.Color = $lastColor

# ---------------------------------------------------------------------
# 2. Update the wording of the "missing key parameter" test title and
#    tweak the column A width / current selection.
# ---------------------------------------------------------------------
$ws.Range("A1:A19").Replace("缺少数据源Id", "name为空")

$ws.Columns("A").ColumnWidth = 56.36

# ---------------------------------------------------------------------
# 3. Rows 3-6: column E switches from 0 (unchecked style) to 1, using
#    the same cell style already applied to E2.
# ---------------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("E3:E9").PasteSpecial(-4122)

$ws.Range("E3:E6").Value = 1

# ---------------------------------------------------------------------
# 4. Add three new rows of test data (rows 7-9), copying the existing
#    row 6 formatting first and then filling in the new content.
# ---------------------------------------------------------------------
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F9").PasteSpecial(-4122)

$ws.Range("A7").Value = "图分析视图-loadData-跨账号不重命名"
$ws.Range("B7").Value = '{"projectId":1334,"graphId":2497,"name":"test"}'
$ws.Range("C7").Value = 401
$ws.Range("D7").Value = "无权访问"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = ""

$ws.Range("A8").Value = "图分析视图-loadData-同账号非当前项目标签不能重命名"
$ws.Range("B8").Value = '{"projectId":1426,"graphId":2827,"name":"test"}'
$ws.Range("C8").Value = 90009
$ws.Range("D8").Value = "无权操作"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = ""

$ws.Range("A9").Value = "图分析视图-loadData-项目下不存在的标签不能重命名"
$ws.Range("B9").Value = '{"projectId":1426,"graphId":2772,"name":"test"}'
$ws.Range("C9").Value = 90009
$ws.Range("D9").Value = "无权操作"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = ""

# ---------------------------------------------------------------------
# 5. Update the active selection to B9, matching the saved view state.
# ---------------------------------------------------------------------
$ws.Range("B9").Select()
